$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") for rows 2 through 28 from serial date 45429 (2024-05-17)
# to serial date 45430 (2024-05-18), keeping existing date formatting.
for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45429) {
        $cell.Value2 = 45430
    }
}
